# Implement first version of lot sizing rules:
# - Bump NrBuckets from 3 to 4 on the Generic sheet.
# - Double the lead-time-like values in Productdata (837 -> 1673) for
#   Part_0001/0002/0003.
# - Add the new (4th) time bucket row to ForecastedAverageDemand and
#   ForcastedStandardDeviation, mirroring the pattern of the existing rows.

$wb = $excel.ActiveWorkbook

# 1. Generic!B4 (NrBuckets): 3 -> 4
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B4").Value = 4

# 2. Productdata!C4:C6: 837 -> 1673
$wsProduct = $wb.Worksheets.Item("Productdata")
$wsProduct.Range("C4").Value = 1673
$wsProduct.Range("C5").Value = 1673
$wsProduct.Range("C6").Value = 1673

# 3. ForecastedAverageDemand: add row 5
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvg.Cells.Item(4, 1).Copy()
$wsAvg.Cells.Item(5, 1).PasteSpecial(-4122)
$wsAvg.Cells.Item(5, 1).Value = 3
$wsAvg.Cells.Item(5, 2).Value = 0
$wsAvg.Cells.Item(5, 3).Value = 0
$wsAvg.Cells.Item(5, 4).Value = 0
$wsAvg.Cells.Item(5, 5).Value = 0
$wsAvg.Cells.Item(5, 6).Value = 0
$wsAvg.Cells.Item(5, 7).Value = 253
$wsAvg.Cells.Item(5, 8).Value = 45
$wsAvg.Cells.Item(5, 9).Value = 75

# 4. ForcastedStandardDeviation: add row 5
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStd.Cells.Item(4, 1).Copy()
$wsStd.Cells.Item(5, 1).PasteSpecial(-4122)
$wsStd.Cells.Item(5, 1).Value = 3
$wsStd.Cells.Item(5, 2).Value = 0
$wsStd.Cells.Item(5, 3).Value = 0
$wsStd.Cells.Item(5, 4).Value = 0
$wsStd.Cells.Item(5, 5).Value = 0
$wsStd.Cells.Item(5, 6).Value = 0
$wsStd.Cells.Item(5, 7).Value = 36.62
$wsStd.Cells.Item(5, 8).Value = 1
$wsStd.Cells.Item(5, 9).Value = 2
